# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2410
#   *_new  -> *_FV2504
# then turn the header row + data range into an Excel Table ("Table1")
# and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells in row 1 (A1:U1) -----------------------
$headerMap = @{
    "Segmentname_old"         = "Segmentname_FV2410"
    "Segmentgruppe_old"       = "Segmentgruppe_FV2410"
    "Segment_old"             = "Segment_FV2410"
    "Datenelement_old"        = "Datenelement_FV2410"
    "Segment ID_old"          = "Segment ID_FV2410"
    "Code_old"                = "Code_FV2410"
    "Qualifier_old"           = "Qualifier_FV2410"
    "Beschreibung_old"        = "Beschreibung_FV2410"
    "Bedingungsausdruck_old"  = "Bedingungsausdruck_FV2410"
    "Bedingung_old"           = "Bedingung_FV2410"
    "Segmentname_new"         = "Segmentname_FV2504"
    "Segmentgruppe_new"       = "Segmentgruppe_FV2504"
    "Segment_new"             = "Segment_FV2504"
    "Datenelement_new"        = "Datenelement_FV2504"
    "Segment ID_new"          = "Segment ID_FV2504"
    "Code_new"                = "Code_FV2504"
    "Qualifier_new"           = "Qualifier_FV2504"
    "Beschreibung_new"        = "Beschreibung_FV2504"
    "Bedingungsausdruck_new"  = "Bedingungsausdruck_FV2504"
    "Bedingung_new"           = "Bedingung_FV2504"
}

$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = [string]$cell.Value2
    if ($headerMap.ContainsKey($old)) {
        $cell.Value = $headerMap[$old]
    }
}

# --- 2. Convert the header + data range into a real Excel Table --------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3. Freeze the header row (row 1) ----------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
